$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Cases" tab query (row 2, column B) previously pulled an extra
# `Cohort` column from the cohort node. That trailing RETURN clause is
# removed so the query only returns up through "Response to Treatment".
# NOTE: use a single-quoted here-string (@' ... '@) so the backticks
# inside the Cypher (column aliases like `Case ID`) are taken literally
# and are NOT treated as PowerShell escape characters.
$newQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
WHERE demo.sex IN ['Male']
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@

$ws.Range("B2").Value = $newQuery

# Reflect the author's final selection/scroll position in the saved view
# (was topLeftCell="A4" / selection B4 -> now top-left default / selection B2).
$ws.Range("B2").Select()
